# Adds the missing "JSON:" textbox to the "detalle-producto" slide,
# listing the JSON payloads used by that screen.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# EMU -> points (PowerPoint COM geometry is expressed in points)
#   off  x=6344066  y=5698965
#   ext  cx=5225097 cy=553998
$left   = 499.53275590551186
$top    = 448.73740157480313
$width  = 411.42496062992126
$height = 43.62188976377953

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = "CuadroTexto 13"

# Match the sibling "CuadroTexto" boxes: no fill, horizontal text,
# word-wrapped, vertically centered, auto-fit to the text.
$shp.Fill.Visible = 0
$tf = $shp.TextFrame
$tf.WordWrap = -1
$tf.AutoSize = 1
$tf.Orientation = 1
$tf.VerticalAnchor = 3
$tf.HorizontalAnchor = 0

$tr = $shp.TextFrame.TextRange
$tr.Text = "JSON:"
$tr.Font.Size = 10
$tr.LanguageID = "es-MX"

$tr2 = $tr.InsertAfter("`rdetalle-producto")
$tr2.Font.Size = 10
$tr2.LanguageID = "es-MX"

$tr3 = $tr.InsertAfter("`rdetalle-producto-recomendados")
$tr3.Font.Size = 10
$tr3.LanguageID = "es-MX"
